# Update the South Fla._B team-specific transition-probability matrix
# after adding more simulated games and speeding up the simulate-game logic.
# Recomputed probabilities (row-normalized frequencies) for the affected
# starting states replace the previous values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2222222222222222
$ws.Range("C2").Value = 0.4962962962962963
$ws.Range("J2").Value = 0.03703703703703703
$ws.Range("P2").Value = 0.1444444444444444
$ws.Range("S2").Value = 0.1
$ws.Range("B3").Value = 0.007092198581560284
$ws.Range("C3").Value = 0.02127659574468085
$ws.Range("J3").Value = 0.02836879432624113
$ws.Range("P3").Value = 0.75177304964539
$ws.Range("S3").Value = 0.1914893617021277
$ws.Range("J4").Value = 0.07547169811320754
$ws.Range("P4").Value = 0.6037735849056604
$ws.Range("S4").Value = 0.3207547169811321
$ws.Range("J5").Value = 0.2
$ws.Range("P5").Value = 0.6
$ws.Range("S5").Value = 0.2
$ws.Range("B6").Value = 0.0851063829787234
$ws.Range("D6").Value = 0.01595744680851064
$ws.Range("F6").Value = 0.02659574468085106
$ws.Range("J6").Value = 0.2606382978723404
$ws.Range("O6").Value = 0.03191489361702127
$ws.Range("Q6").Value = 0.1223404255319149
$ws.Range("R6").Value = 0.05319148936170213
$ws.Range("S6").Value = 0.4042553191489361
$ws.Range("B7").Value = 0.1237113402061856
$ws.Range("D7").Value = 0.02577319587628866
$ws.Range("F7").Value = 0.05670103092783505
$ws.Range("J7").Value = 0.1185567010309278
$ws.Range("O7").Value = 0.02061855670103093
$ws.Range("Q7").Value = 0.1701030927835052
$ws.Range("R7").Value = 0.05670103092783505
$ws.Range("S7").Value = 0.4278350515463917
$ws.Range("B8").Value = 0.08894878706199461
$ws.Range("D8").Value = 0.01886792452830189
$ws.Range("E8").Value = 0.002695417789757413
$ws.Range("F8").Value = 0.04851752021563342
$ws.Range("J8").Value = 0.1347708894878706
$ws.Range("O8").Value = 0.01617250673854448
$ws.Range("Q8").Value = 0.1967654986522911
$ws.Range("R8").Value = 0.07547169811320754
$ws.Range("S8").Value = 0.4177897574123989
$ws.Range("B9").Value = 0.07582938388625593
$ws.Range("D9").Value = 0.01421800947867299
$ws.Range("F9").Value = 0.1090047393364929
$ws.Range("J9").Value = 0.1327014218009479
$ws.Range("O9").Value = 0.02369668246445497
$ws.Range("Q9").Value = 0.1374407582938389
$ws.Range("R9").Value = 0.09004739336492891
$ws.Range("S9").Value = 0.4170616113744076
$ws.Range("B10").Value = 0.1053089643167972
$ws.Range("D10").Value = 0.03220191470844212
$ws.Range("E10").Value = 0.006092254134029591
$ws.Range("F10").Value = 0.07049608355091384
$ws.Range("J10").Value = 0.1105308964316797
$ws.Range("O10").Value = 0.01740644038294169
$ws.Range("Q10").Value = 0.1836379460400348
$ws.Range("R10").Value = 0.08006962576153177
$ws.Range("S10").Value = 0.3942558746736293
$ws.Range("G11").Value = 0.1548387096774194
$ws.Range("J11").Value = 0.08064516129032258
$ws.Range("K11").Value = 0.1967741935483871
$ws.Range("L11").Value = 0.5580645161290323
$ws.Range("S11").Value = 0.00967741935483871
$ws.Range("G12").Value = 0.7955801104972375
$ws.Range("J12").Value = 0.1491712707182321
$ws.Range("K12").Value = 0.01104972375690608
$ws.Range("L12").Value = 0.03867403314917127
$ws.Range("S12").Value = 0.005524861878453038
$ws.Range("G13").Value = 0.4146341463414634
$ws.Range("J13").Value = 0.5365853658536586
$ws.Range("S13").Value = 0.04878048780487805
$ws.Range("J14").Value = 1
$ws.Range("F15").Value = 0.015625
$ws.Range("H15").Value = 0.109375
$ws.Range("I15").Value = 0.09895833333333333
$ws.Range("J15").Value = 0.359375
$ws.Range("K15").Value = 0.07291666666666667
$ws.Range("M15").Value = 0.005208333333333333
$ws.Range("O15").Value = 0.046875
$ws.Range("S15").Value = 0.2916666666666667
$ws.Range("H16").Value = 0.1686046511627907
$ws.Range("I16").Value = 0.08139534883720931
$ws.Range("J16").Value = 0.3662790697674418
$ws.Range("K16").Value = 0.1220930232558139
$ws.Range("M16").Value = 0.01744186046511628
$ws.Range("O16").Value = 0.04651162790697674
$ws.Range("S16").Value = 0.1976744186046512
$ws.Range("F17").Value = 0.00819672131147541
$ws.Range("H17").Value = 0.1639344262295082
$ws.Range("I17").Value = 0.1284153005464481
$ws.Range("J17").Value = 0.4207650273224044
$ws.Range("K17").Value = 0.1010928961748634
$ws.Range("M17").Value = 0.01092896174863388
$ws.Range("O17").Value = 0.04918032786885246
$ws.Range("S17").Value = 0.1174863387978142
$ws.Range("F18").Value = 0.025
$ws.Range("H18").Value = 0.175
$ws.Range("I18").Value = 0.11875
$ws.Range("J18").Value = 0.44375
$ws.Range("K18").Value = 0.08749999999999999
$ws.Range("M18").Value = 0.0125
$ws.Range("O18").Value = 0.05
$ws.Range("S18").Value = 0.08749999999999999
$ws.Range("F19").Value = 0.009083402146985962
$ws.Range("H19").Value = 0.1907514450867052
$ws.Range("I19").Value = 0.0916597853014038
$ws.Range("J19").Value = 0.365813377374071
$ws.Range("K19").Value = 0.129644921552436
$ws.Range("M19").Value = 0.02559867877786953
$ws.Range("N19").Value = 0.0008257638315441783
$ws.Range("O19").Value = 0.07101568951279934
$ws.Range("S19").Value = 0.115606936416185
